$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ P = "320017958830"; Q = "`$19.04";    R = "PASS" }
    3  = @{ P = "320017958841"; Q = "`$27.50";    R = "PASS" }
    4  = @{ P = "320017958874"; Q = "`$31.73";    R = "PASS" }
    5  = @{ P = "320017958896"; Q = "`$43.36";    R = "PASS" }
    6  = @{ P = "320017958933"; Q = "`$56.05";    R = "PASS" }
    7  = @{ P = "320017958955"; Q = "`$231.08";   R = "PASS" }
    8  = @{ P = "320017958988"; Q = "`$19.04";    R = "PASS" }
    9  = @{ P = "320017959002"; Q = "`$23.27";    R = "PASS" }
    10 = @{ P = "320017959035"; Q = "`$27.50";    R = "PASS" }
    11 = @{ P = "320017959057"; Q = "`$40.19";    R = "PASS" }
    12 = @{ P = "320017959090"; Q = "`$52.88";    R = "PASS" }
    13 = @{ P = "320017959116"; Q = "`$14.81";    R = "PASS" }
    14 = @{ P = "320017959149"; Q = "`$17.98";    R = "PASS" }
    15 = @{ P = "320017959160"; Q = "`$21.15";    R = "PASS" }
    16 = @{ P = "320017959208"; Q = "`$31.73";    R = "PASS" }
    17 = @{ P = "320017959220"; Q = "`$42.30";    R = "PASS" }
    18 = @{ P = "320017959263"; Q = "`$85.66";    R = "FAIL" }
    19 = @{ P = "320017959285"; Q = "`$53.93";    R = "PASS" }
    20 = @{ P = "320017959311"; Q = "`$85.66";    R = "FAIL" }
    21 = @{ P = "320017959333"; Q = "`$111.04";   R = "PASS" }
    22 = @{ P = "320017959366"; Q = "`$233.07";   R = "FAIL" }
    23 = @{ P = "320017959377"; Q = "`$476.72";   R = "FAIL" }
    24 = @{ P = "320017959388"; Q = "`$306.68";   R = "FAIL" }
    25 = @{ P = "320017959403"; Q = "`$52.88";    R = "PASS" }
    26 = @{ P = "320017959414"; Q = "`$1,171.41"; R = "PASS" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $pCell = $ws.Range("P$row")
    $pCell.NumberFormat = "@"
    $pCell.Value = $vals.P
    $pCell.ClearFormats()

    $qCell = $ws.Range("Q$row")
    $qCell.NumberFormat = "@"
    $qCell.Value = $vals.Q
    $qCell.ClearFormats()

    $ws.Range("R$row").Value = $vals.R
}
